$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1154.7
$ws.Range("I4").Value = 864
$ws.Range("J4").Value = 1833
$ws.Range("K4").Value = 864
$ws.Range("L4").Value = 1833
$ws.Range("M4").Value = -750
$ws.Range("N4").Value = -2061

$ws.Range("H39").Value = 2537.4443
$ws.Range("I39").Value = 334
$ws.Range("K39").Value = 1002
$ws.Range("M39").Value = -706

$ws.Range("H40").Value = 2616.5
$ws.Range("J40").Value = 3061.25
$ws.Range("L40").Value = 3061.25
$ws.Range("N40").Value = -3411.25

$ws.Range("H70").Value = 9067.9375
$ws.Range("I70").Value = 6329.3335
$ws.Range("J70").Value = 9699.923
$ws.Range("K70").Value = 18988.0005
$ws.Range("L70").Value = 29099.769
$ws.Range("M70").Value = -18718.0005
$ws.Range("N70").Value = -29639.769

$ws.Range("H73").Value = 9067.9375
$ws.Range("I73").Value = 6329.3335
$ws.Range("J73").Value = 9699.923
$ws.Range("K73").Value = 18988.0005
$ws.Range("L73").Value = 29099.769
$ws.Range("M73").Value = -18052.0005
$ws.Range("N73").Value = -30971.769

$ws.Range("H98").Value = 2923.738
$ws.Range("I98").Value = 2475.0588
$ws.Range("J98").Value = 4830.625
$ws.Range("K98").Value = 2475.0588
$ws.Range("L98").Value = 4830.625
$ws.Range("M98").Value = -977.0587999999998
$ws.Range("N98").Value = -7826.625

$ws.Range("H122").Value = 2923.738
$ws.Range("I122").Value = 2475.0588
$ws.Range("J122").Value = 4830.625
$ws.Range("K122").Value = 7425.176399999999
$ws.Range("L122").Value = 14491.875
$ws.Range("M122").Value = -4975.176399999999
$ws.Range("N122").Value = -19391.875

$ws.Range("H131").Value = 3268.9167
$ws.Range("J131").Value = 9360
$ws.Range("L131").Value = 28080
$ws.Range("N131").Value = -38160

$ws.Range("H132").Value = 3386.2144
$ws.Range("I132").Value = 3357.6458
$ws.Range("K132").Value = 10072.9374
$ws.Range("M132").Value = -7542.937399999999

$ws.Range("H138").Value = 3076.0488
$ws.Range("J138").Value = 2653.7144
$ws.Range("L138").Value = 7961.1432
$ws.Range("N138").Value = -18241.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H74").Value = 4446.909
$ws.Range("I74").Value = 3333.75
$ws.Range("K74").Value = 3333.75
$ws.Range("M74").Value = -2459.75

$ws.Range("H77").Value = 4446.909
$ws.Range("I77").Value = 3333.75
$ws.Range("K77").Value = 16668.75
$ws.Range("M77").Value = -12300.75

$ws.Range("H97").Value = 4359.5884
$ws.Range("J97").Value = 2149.6
$ws.Range("L97").Value = 2149.6
$ws.Range("N97").Value = -3141.6

$ws.Range("H102").Value = 1268.5122
$ws.Range("I102").Value = 1186.4324
$ws.Range("K102").Value = 1186.4324
$ws.Range("M102").Value = 435.5676000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 41689
$ws.Range("I82").Value = 9483.333
$ws.Range("J82").Value = 89997.5
$ws.Range("K82").Value = 9483.333
$ws.Range("L82").Value = 89997.5
$ws.Range("M82").Value = -9100.333
$ws.Range("N82").Value = -90763.5

$ws.Range("H85").Value = 41689
$ws.Range("I85").Value = 9483.333
$ws.Range("J85").Value = 89997.5
$ws.Range("K85").Value = 9483.333
$ws.Range("L85").Value = 89997.5
$ws.Range("M85").Value = -8157.333000000001
$ws.Range("N85").Value = -92649.5

$ws.Range("H94").Value = 5836.476
$ws.Range("I94").Value = 510.8125
$ws.Range("K94").Value = 510.8125
$ws.Range("M94").Value = -59.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1650
$ws.Range("I23").Value = 1650
$ws.Range("K23").Value = 1650
$ws.Range("M23").Value = -1410

$ws.Range("H27").Value = 1650
$ws.Range("I27").Value = 1650
$ws.Range("K27").Value = 1650
$ws.Range("M27").Value = -1458

$ws.Range("H58").Value = 3772.7334
$ws.Range("I58").Value = 2017.8
$ws.Range("J58").Value = 7282.6
$ws.Range("K58").Value = 2017.8
$ws.Range("L58").Value = 7282.6
$ws.Range("M58").Value = -1814.8
$ws.Range("N58").Value = -7688.6

$ws.Range("H136").Value = 3772.7334
$ws.Range("I136").Value = 2017.8
$ws.Range("J136").Value = 7282.6
$ws.Range("K136").Value = 6053.4
$ws.Range("L136").Value = 21847.8
$ws.Range("M136").Value = -3503.4
$ws.Range("N136").Value = -26947.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 539610.94
$ws.Range("I4").Value = 902503.5
$ws.Range("J4").Value = 104139.8
$ws.Range("K4").Value = 2707510.5
$ws.Range("L4").Value = 312419.4
$ws.Range("M4").Value = -2707398.5
$ws.Range("N4").Value = -312643.4

$ws.Range("H19").Value = 1499.5
$ws.Range("J19").Value = 1999
$ws.Range("L19").Value = 5997
$ws.Range("N19").Value = -6345

$ws.Range("H107").Value = 1066.4117
$ws.Range("I107").Value = 649.7143
$ws.Range("J107").Value = 1358.1
$ws.Range("K107").Value = 1949.1429
$ws.Range("L107").Value = 4074.3
$ws.Range("M107").Value = -29.14289999999983
$ws.Range("N107").Value = -7914.299999999999

$ws.Range("H113").Value = 1130.0769
$ws.Range("J113").Value = 1088.75
$ws.Range("L113").Value = 3266.25
$ws.Range("N113").Value = -7606.25

$ws.Range("H136").Value = 20238
$ws.Range("I136").Value = 7243.1665
$ws.Range("J136").Value = 29984.125
$ws.Range("K136").Value = 21729.4995
$ws.Range("L136").Value = 89952.375
$ws.Range("M136").Value = -16629.4995
$ws.Range("N136").Value = -100152.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4999
$ws.Range("I113").Value = 4999
$ws.Range("K113").Value = 4999
$ws.Range("M113").Value = -2829

$ws.Range("H134").Value = 84679.8
$ws.Range("J134").Value = 84679.8
$ws.Range("L134").Value = 254039.4
$ws.Range("N134").Value = -259109.4

$ws.Range("H136").Value = 35496.777
$ws.Range("J136").Value = 35496.777
$ws.Range("L136").Value = 106490.331
$ws.Range("N136").Value = -111590.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1689.3572
$ws.Range("J46").Value = 2599.3333
$ws.Range("L46").Value = 2599.3333
$ws.Range("N46").Value = -2975.3333

$ws.Range("H122").Value = 6392.1875
$ws.Range("J122").Value = 6685
$ws.Range("L122").Value = 20055
$ws.Range("N122").Value = -24955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H74").Value = 26277.5
$ws.Range("I74").Value = 22555
$ws.Range("J74").Value = 30000
$ws.Range("K74").Value = 22555
$ws.Range("L74").Value = 30000
$ws.Range("M74").Value = -21619
$ws.Range("N74").Value = -31872

$ws.Range("H77").Value = 26277.5
$ws.Range("I77").Value = 22555
$ws.Range("J77").Value = 30000
$ws.Range("K77").Value = 67665
$ws.Range("L77").Value = 90000
$ws.Range("M77").Value = -62985
$ws.Range("N77").Value = -99360

$ws.Range("H81").Value = 3732.04
$ws.Range("I81").Value = 4417.737
$ws.Range("J81").Value = 1560.6666
$ws.Range("K81").Value = 8835.474
$ws.Range("L81").Value = 3121.3332
$ws.Range("M81").Value = -7774.474
$ws.Range("N81").Value = -5243.3332

$ws.Range("H84").Value = 3732.04
$ws.Range("I84").Value = 4417.737
$ws.Range("J84").Value = 1560.6666
$ws.Range("K84").Value = 44177.37
$ws.Range("L84").Value = 15606.666
$ws.Range("M84").Value = -38873.37
$ws.Range("N84").Value = -26214.666

